# Generate Report for Handback
#
# Row 7 of both the "zh-cn" and "de-de" sheets represents the handback
# status for 1f07d7dd-d537-48ae-99fb-fc95bbbbe555.md. A new handback was
# received, but it was generated from a stale source revision, so the
# report records the version-mismatch error instead of a clean handback.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58ca70f0813117e924f309f8d7fae10efad26b85/e2e/1f07d7dd-d537-48ae-99fb-fc95bbbbe555.md"
$handbackDisplay = "1f07d7dd-d537-48ae-99fb-fc95bbbbe555.md"

function Update-HandbackRow7 {
    param(
        $ws,
        [string]$targetXlf,
        [string]$handbackDateTime,
        [string]$errorDetail
    )

    # Latest Target File
    $ws.Range("J7").Value = $targetXlf
    # Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime
    # Error Detail
    $ws.Range("P7").Value = $errorDetail

    # Latest Handback File - text + hyperlink (same workflow as column A's
    # handoff-file links and the existing I2:I5 handback-file links).
    $ws.Hyperlinks.Add($ws.Range("I7"), $handbackUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $handbackDisplay) | Out-Null
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow7 -ws $wsZhCn `
    -targetXlf "1f07d7dd-d537-48ae-99fb-fc95bbbbe555.91c4bfbd773bee313cfd88079d6017730dbb5b8e.zh-cn.xlf" `
    -handbackDateTime "2016-09-04 11:00:58" `
    -errorDetail "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54a4a9e009aa110ded4db3ac9ec09514253901f5/e2e/1f07d7dd-d537-48ae-99fb-fc95bbbbe555.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58ca70f0813117e924f309f8d7fae10efad26b85/e2e/1f07d7dd-d537-48ae-99fb-fc95bbbbe555.md."

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow7 -ws $wsDeDe `
    -targetXlf "1f07d7dd-d537-48ae-99fb-fc95bbbbe555.91c4bfbd773bee313cfd88079d6017730dbb5b8e.de-de.xlf" `
    -handbackDateTime "2016-09-04 11:01:11" `
    -errorDetail "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54a4a9e009aa110ded4db3ac9ec09514253901f5/e2e/1f07d7dd-d537-48ae-99fb-fc95bbbbe555.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58ca70f0813117e924f309f8d7fae10efad26b85/e2e/1f07d7dd-d537-48ae-99fb-fc95bbbbe555.md."
